$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22; this shifts existing rows 22-111 down to 23-112
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with the new daily price record.
# Non-changing columns are copied from the (now shifted) row 23, which held
# the same Mercado / Producto metadata.
$ws.Range("A22").Value = 10
$ws.Range("B22").Value = "Vega Modelo de Temuco"
$ws.Range("C22").Value = "La Araucanía"
$ws.Range("D22").Value = 45145
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100108
$ws.Range("H22").Value = "Tropicales y subtropicales"
$ws.Range("I22").Value = 100108003
$ws.Range("J22").Value = "Maracuyá"
$ws.Range("K22").Value = "Sin especificar"
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 40
$ws.Range("N22").Value = 42000
$ws.Range("O22").Value = 42000
$ws.Range("P22").Value = 42000
$ws.Range("Q22").Value = "$/caja 18 kilos"
$ws.Range("R22").Value = "Región de Arica y Parinacota"
$ws.Range("S22").Value = 2333
$ws.Range("T22").Value = 18
